# Convert the numeric month values in column C (rows 5-84) into their
# Spanish three-letter abbreviations (e.g. 8 -> "Ago.", 7 -> "Jul.", ...).
# This mirrors the author's change of turning the "Mes" column from a
# plain integer into a text label, backed by shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 5; $row -le 84; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C = "Mes"
    $monthNumber = [int]$cell.Value2
    $cell.Value = $monthNames[$monthNumber]
}

# Also move the workbook window position slightly, matching the recorded
# view-state change in the workbook part.
$excel.ActiveWindow.Left = 23880
